$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 (duplicate fastq file entry:
# "Brent_3275_Small_20-1_GTAC_1_SIC_Index2_07_TGAGGTT_GAGTTGAG_S31_R1_001.fastq.gz")
# This shifts all subsequent rows up by one.
$ws.Rows.Item(5).Delete()

# Select the new row 5 (mirrors the workbook's saved selection state after the edit)
$ws.Range("A5:XFD5").Select()
